$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 25003392
$ws.Range("I106").Value = 28574782
$ws.Range("K106").Value = 28574782
$ws.Range("M106").Value = -28574151

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1809.7059
$ws.Range("I111").Value = 2001.5
$ws.Range("J111").Value = 1535.7142
$ws.Range("K111").Value = 6004.5
$ws.Range("L111").Value = 4607.142599999999
$ws.Range("M111").Value = -2937.5
$ws.Range("N111").Value = -10741.1426

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1853.96
$ws.Range("I138").Value = 861.55884
$ws.Range("J138").Value = 2365.197
$ws.Range("K138").Value = 2584.67652
$ws.Range("L138").Value = 7095.591
$ws.Range("M138").Value = 2555.32348
$ws.Range("N138").Value = -17375.591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 22085.25
$ws.Range("I21").Value = 4610.8335
$ws.Range("J21").Value = 74508.5
$ws.Range("K21").Value = 4610.8335
$ws.Range("L21").Value = 74508.5
$ws.Range("M21").Value = -4236.8335
$ws.Range("N21").Value = -75256.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1758.159
$ws.Range("I74").Value = 998.15
$ws.Range("K74").Value = 998.15
$ws.Range("M74").Value = -124.15

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1758.159
$ws.Range("I77").Value = 998.15
$ws.Range("K77").Value = 4990.75
$ws.Range("M77").Value = -622.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 78201.14
$ws.Range("J101").Value = 78201.14
$ws.Range("L101").Value = 78201.14
$ws.Range("N101").Value = -84691.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2207.8545
$ws.Range("I132").Value = 1645.2094
$ws.Range("K132").Value = 4935.6282
$ws.Range("M132").Value = -2405.6282

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 595.1111
$ws.Range("I64").Value = 671.2
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 671.2
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = -446.2
$ws.Range("N64").Value = -950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 595.1111
$ws.Range("I67").Value = 671.2
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 671.2
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 108.8
$ws.Range("N67").Value = -2060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2362.8572
$ws.Range("I86").Value = 2855.2307
$ws.Range("J86").Value = 1562.75
$ws.Range("K86").Value = 2855.2307
$ws.Range("L86").Value = 1562.75
$ws.Range("M86").Value = -1732.2307
$ws.Range("N86").Value = -3808.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2362.8572
$ws.Range("I89").Value = 2855.2307
$ws.Range("J89").Value = 1562.75
$ws.Range("K89").Value = 14276.1535
$ws.Range("L89").Value = 7813.75
$ws.Range("M89").Value = -8660.1535
$ws.Range("N89").Value = -19045.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1455.6364
$ws.Range("J99").Value = 1960.4
$ws.Range("L99").Value = 1960.4
$ws.Range("N99").Value = -4956.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8929862
$ws.Range("I105").Value = 10417949
$ws.Range("J105").Value = 1337.5
$ws.Range("K105").Value = 10417949
$ws.Range("L105").Value = 1337.5
$ws.Range("M105").Value = -10416202
$ws.Range("N105").Value = -4831.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 2750
$ws.Range("I38").Value = 1500
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 1500
$ws.Range("L38").Value = 4000
$ws.Range("M38").Value = -1123
$ws.Range("N38").Value = -4754

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H46").Value = 2750
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -1289
$ws.Range("N46").Value = -4422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1538.125
$ws.Range("I58").Value = 1184.3334
$ws.Range("J58").Value = 2127.7778
$ws.Range("K58").Value = 1184.3334
$ws.Range("L58").Value = 2127.7778
$ws.Range("M58").Value = -981.3334
$ws.Range("N58").Value = -2533.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 911.61536
$ws.Range("I105").Value = 879.55
$ws.Range("J105").Value = 1018.5
$ws.Range("K105").Value = 879.55
$ws.Range("L105").Value = 1018.5
$ws.Range("M105").Value = 867.45
$ws.Range("N105").Value = -4512.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4903782
$ws.Range("I132").Value = 1628.65
$ws.Range("J132").Value = 11906858
$ws.Range("K132").Value = 4885.950000000001
$ws.Range("L132").Value = 35720574
$ws.Range("M132").Value = -2355.950000000001
$ws.Range("N132").Value = -35725634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1538.125
$ws.Range("I136").Value = 1184.3334
$ws.Range("J136").Value = 2127.7778
$ws.Range("K136").Value = 3553.0002
$ws.Range("L136").Value = 6383.3334
$ws.Range("M136").Value = -1003.0002
$ws.Range("N136").Value = -11483.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 309.0909
$ws.Range("I7").Value = 235.71428
$ws.Range("J7").Value = 437.5
$ws.Range("K7").Value = 707.14284
$ws.Range("L7").Value = 1312.5
$ws.Range("M7").Value = -595.14284
$ws.Range("N7").Value = -1536.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51.347828
$ws.Range("I12").Value = 43.833332
$ws.Range("J12").Value = 54
$ws.Range("K12").Value = 131.499996
$ws.Range("L12").Value = 162
$ws.Range("M12").Value = 41.50000399999999
$ws.Range("N12").Value = -508

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 414.51724
$ws.Range("I15").Value = 74
$ws.Range("J15").Value = 453.80768
$ws.Range("K15").Value = 222
$ws.Range("L15").Value = 1361.42304
$ws.Range("M15").Value = -82
$ws.Range("N15").Value = -1641.42304

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1249.6666
$ws.Range("I17").Value = 375
$ws.Range("J17").Value = 2999
$ws.Range("K17").Value = 1125
$ws.Range("L17").Value = 8997
$ws.Range("M17").Value = -956
$ws.Range("N17").Value = -9335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 992.1053000000001
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 991.17645
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 2973.52935
$ws.Range("M20").Value = -2773
$ws.Range("N20").Value = -3427.52935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1025
$ws.Range("I22").Value = 733.3333
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 2199.9999
$ws.Range("L22").Value = 3600
$ws.Range("M22").Value = -2030.9999
$ws.Range("N22").Value = -3938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 596.875
$ws.Range("I25").Value = 83.333336
$ws.Range("J25").Value = 715.38464
$ws.Range("K25").Value = 250.000008
$ws.Range("L25").Value = 2146.15392
$ws.Range("M25").Value = -81.00000800000001
$ws.Range("N25").Value = -2484.15392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1025
$ws.Range("I27").Value = 733.3333
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 2199.9999
$ws.Range("L27").Value = 3600
$ws.Range("M27").Value = -2097.9999
$ws.Range("N27").Value = -3804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 41667240
$ws.Range("I29").Value = 200.5
$ws.Range("J29").Value = 55556252
$ws.Range("K29").Value = 601.5
$ws.Range("L29").Value = 166668756
$ws.Range("M29").Value = -324.5
$ws.Range("N29").Value = -166669310

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 596.875
$ws.Range("I30").Value = 83.333336
$ws.Range("J30").Value = 715.38464
$ws.Range("K30").Value = 250.000008
$ws.Range("L30").Value = 2146.15392
$ws.Range("M30").Value = -148.000008
$ws.Range("N30").Value = -2350.15392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4000
$ws.Range("I56").Value = 4000
$ws.Range("K56").Value = 4000
$ws.Range("M56").Value = -3470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1454
$ws.Range("I63").Value = 934
$ws.Range("J63").Value = 3014
$ws.Range("K63").Value = 2802
$ws.Range("L63").Value = 9042
$ws.Range("M63").Value = -2053
$ws.Range("N63").Value = -10540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 1454
$ws.Range("I66").Value = 934
$ws.Range("J66").Value = 3014
$ws.Range("K66").Value = 8406
$ws.Range("L66").Value = 27126
$ws.Range("M66").Value = -4662
$ws.Range("N66").Value = -34614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 248
$ws.Range("I75").Value = 203.33333
$ws.Range("J75").Value = 315
$ws.Range("K75").Value = 609.99999
$ws.Range("L75").Value = 945
$ws.Range("M75").Value = 388.00001
$ws.Range("N75").Value = -2941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 248
$ws.Range("I78").Value = 203.33333
$ws.Range("J78").Value = 315
$ws.Range("K78").Value = 1829.99997
$ws.Range("L78").Value = 2835
$ws.Range("M78").Value = 3162.00003
$ws.Range("N78").Value = -12819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 22222504
$ws.Range("J107").Value = 71428830
$ws.Range("L107").Value = 214286490
$ws.Range("N107").Value = -214290330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 9000
$ws.Range("I26").Value = 9000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -8705
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2766.8262
$ws.Range("I61").Value = 2425.4707
$ws.Range("J61").Value = 3734
$ws.Range("K61").Value = 2425.4707
$ws.Range("L61").Value = 3734
$ws.Range("M61").Value = -2223.4707
$ws.Range("N61").Value = -4138

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2766.8262
$ws.Range("I113").Value = 2425.4707
$ws.Range("J113").Value = 3734
$ws.Range("K113").Value = 2425.4707
$ws.Range("L113").Value = 3734
$ws.Range("M113").Value = -255.4706999999999
$ws.Range("N113").Value = -8074

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 80000
$ws.Range("J63").Value = 80000
$ws.Range("L63").Value = 80000
$ws.Range("N63").Value = -81248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 80000
$ws.Range("J66").Value = 80000
$ws.Range("L66").Value = 240000
$ws.Range("N66").Value = -246240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 536.875
$ws.Range("I100").Value = 358.4
$ws.Range("J100").Value = 834.3333
$ws.Range("K100").Value = 716.8
$ws.Range("L100").Value = 1668.6666
$ws.Range("M100").Value = -175.8
$ws.Range("N100").Value = -2750.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 59333.332
$ws.Range("J118").Value = 59333.332
$ws.Range("L118").Value = 59333.332
$ws.Range("N118").Value = -62647.332
